# yearly_death_totals.xlsx — add back the 2010-2012 columns in front of the
# existing 2013-2020 series (commit: "Added back in 2010-2012 because why not").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:H data three columns to the right, opening up A:C.
$ws.Range("A:C").Insert(-4161)   # xlShiftToRight

# New year headers in row 1.
$ws.Range("A1").Value = 2010
$ws.Range("B1").Value = 2011
$ws.Range("C1").Value = 2012

# New yearly death totals in row 2.
$ws.Range("A2").Value = 19983
$ws.Range("B2").Value = 21414
$ws.Range("C2").Value = 21837

# The author's active cell ended up on the first cell of the original data
# (now shifted to D2).
$ws.Range("D2").Select()

# Best-effort: the saved workbook window was resized/repositioned too.
$win = $wb.Windows.Item(1)
$win.Left = 11250
$win.Top = 0
$win.Width = 11250
$win.Height = 15000
